$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 56
$ws.Range("K4").Value = 56
$ws.Range("M4").Value = 58

$ws.Range("H17").Value = 555591.0600000001
$ws.Range("J17").Value = 611113.4
$ws.Range("L17").Value = 1833340.2
$ws.Range("N17").Value = -1833676.2

$ws.Range("H18").Value = 817.5
$ws.Range("I18").Value = 812.8570999999999
$ws.Range("J18").Value = 850
$ws.Range("K18").Value = 812.8570999999999
$ws.Range("L18").Value = 850
$ws.Range("M18").Value = -528.8570999999999
$ws.Range("N18").Value = -1418

$ws.Range("H31").Value = 1491.5
$ws.Range("I31").Value = 1189.4
$ws.Range("K31").Value = 3568.2
$ws.Range("M31").Value = -3338.2

$ws.Range("H39").Value = 2537.2354
$ws.Range("I39").Value = 1079.25
$ws.Range("K39").Value = 3237.75
$ws.Range("M39").Value = -2941.75

$ws.Range("H40").Value = 4519.9
$ws.Range("I40").Value = 1899.875
$ws.Range("K40").Value = 1899.875
$ws.Range("M40").Value = -1724.875

$ws.Range("H98").Value = 1403.9166
$ws.Range("I98").Value = 1030.174
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 1030.174
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = 467.826
$ws.Range("N98").Value = -12996

$ws.Range("H112").Value = 102106.9
$ws.Range("J112").Value = 69132.60000000001
$ws.Range("L112").Value = 207397.8
$ws.Range("N112").Value = -209613.8

$ws.Range("H122").Value = 1403.9166
$ws.Range("I122").Value = 1030.174
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 3090.522
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -640.5219999999999
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 2695.1
$ws.Range("I132").Value = 2679.0527
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8037.158100000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5507.158100000001
$ws.Range("N132").Value = -14060

$ws.Range("H137").Value = 3222
$ws.Range("I137").Value = 2983.36
$ws.Range("K137").Value = 8950.08
$ws.Range("M137").Value = -6400.08

$ws.Range("H138").Value = 2637.9
$ws.Range("I138").Value = 1785
$ws.Range("K138").Value = 5355
$ws.Range("M138").Value = -215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3415.4546
$ws.Range("I32").Value = 1900.4509
$ws.Range("J32").Value = 22731.75
$ws.Range("K32").Value = 1900.4509
$ws.Range("L32").Value = 22731.75
$ws.Range("M32").Value = -1613.4509
$ws.Range("N32").Value = -23305.75

$ws.Range("H122").Value = 4938.9375
$ws.Range("I122").Value = 4367.154
$ws.Range("J122").Value = 7416.6665
$ws.Range("K122").Value = 13101.462
$ws.Range("L122").Value = 22249.9995
$ws.Range("M122").Value = -10651.462
$ws.Range("N122").Value = -27149.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 254.91667
$ws.Range("I7").Value = 306.55554
$ws.Range("K7").Value = 306.55554
$ws.Range("M7").Value = -193.55554

$ws.Range("H11").Value = 599
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H19").Value = 1500.4
$ws.Range("J19").Value = 1800
$ws.Range("L19").Value = 1800
$ws.Range("N19").Value = -2140

$ws.Range("H24").Value = 1500.4
$ws.Range("J24").Value = 1800
$ws.Range("L24").Value = 1800
$ws.Range("N24").Value = -2140

$ws.Range("H31").Value = 2355.9443
$ws.Range("I31").Value = 2714.4614
$ws.Range("J31").Value = 1423.8
$ws.Range("K31").Value = 2714.4614
$ws.Range("L31").Value = 1423.8
$ws.Range("M31").Value = -2419.4614
$ws.Range("N31").Value = -2013.8

$ws.Range("H34").Value = 2355.9443
$ws.Range("I34").Value = 2714.4614
$ws.Range("J34").Value = 1423.8
$ws.Range("K34").Value = 2714.4614
$ws.Range("L34").Value = 1423.8
$ws.Range("M34").Value = -2512.4614
$ws.Range("N34").Value = -1827.8

$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -9407
$ws.Range("N42").ClearContents()

$ws.Range("H132").Value = 62501956
$ws.Range("I132").Value = 66668680
$ws.Range("K132").Value = 200006040
$ws.Range("M132").Value = -200003510

$ws.Range("H134").Value = 6788934
$ws.Range("I134").Value = 7611176.5
$ws.Range("K134").Value = 22833529.5
$ws.Range("M134").Value = -22830994.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 278
$ws.Range("I6").Value = 257
$ws.Range("K6").Value = 771
$ws.Range("M6").Value = -658

$ws.Range("H10").Value = 176
$ws.Range("I10").Value = 176
$ws.Range("K10").Value = 528
$ws.Range("M10").Value = -389

$ws.Range("H37").Value = 115969.39
$ws.Range("J37").Value = 115969.39
$ws.Range("L37").Value = 347908.17
$ws.Range("N37").Value = -348132.17

$ws.Range("H59").Value = 3324.75
$ws.Range("J59").Value = 2300
$ws.Range("L59").Value = 6900
$ws.Range("N59").Value = -7980

$ws.Range("H75").Value = 333.75
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 333.75
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 1001.25
$ws.Range("N75").Value = -2997.25
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 333.75
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 333.75
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 3003.75
$ws.Range("N78").Value = -12987.75
$ws.Range("M78").ClearContents()

$ws.Range("H108").Value = 3624.5
$ws.Range("I108").Value = 3250
$ws.Range("K108").Value = 9750
$ws.Range("M108").Value = -6870

$ws.Range("H113").Value = 143986
$ws.Range("J113").Value = 1250.25
$ws.Range("L113").Value = 3750.75
$ws.Range("N113").Value = -8090.75

$ws.Range("H121").Value = 127531.375
$ws.Range("I121").Value = 200989.8
$ws.Range("J121").Value = 5100.6665
$ws.Range("K121").Value = 602969.3999999999
$ws.Range("L121").Value = 15301.9995
$ws.Range("M121").Value = -601659.3999999999
$ws.Range("N121").Value = -17921.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2061.8
$ws.Range("I80").Value = 2103
$ws.Range("K80").Value = 2103
$ws.Range("M80").Value = -1105

$ws.Range("H83").Value = 2061.8
$ws.Range("I83").Value = 2103
$ws.Range("K83").Value = 10515
$ws.Range("M83").Value = -5523

$ws.Range("H132").Value = 5004210
$ws.Range("I132").Value = 6253980
$ws.Range("J132").Value = 5130
$ws.Range("K132").Value = 18761940
$ws.Range("L132").Value = 15390
$ws.Range("M132").Value = -18759410
$ws.Range("N132").Value = -20450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4297
$ws.Range("I7").Value = 4168.5713
$ws.Range("J7").Value = 4596.6665
$ws.Range("K7").Value = 4168.5713
$ws.Range("L7").Value = 4596.6665
$ws.Range("M7").Value = -4056.5713
$ws.Range("N7").Value = -4820.6665

$ws.Range("H16").Value = 1570.5454
$ws.Range("I16").Value = 1194
$ws.Range("J16").Value = 2229.5
$ws.Range("K16").Value = 1194
$ws.Range("L16").Value = 2229.5
$ws.Range("M16").Value = -1024
$ws.Range("N16").Value = -2569.5

$ws.Range("H55").Value = 630.86664
$ws.Range("I55").Value = 279.14285
$ws.Range("K55").Value = 279.14285
$ws.Range("M55").Value = -106.14285

$ws.Range("H61").Value = 3748.15
$ws.Range("I61").Value = 3798
$ws.Range("K61").Value = 3798
$ws.Range("M61").Value = -3596

$ws.Range("H113").Value = 3748.15
$ws.Range("I113").Value = 3798
$ws.Range("K113").Value = 3798
$ws.Range("M113").Value = -1628

$ws.Range("H126").Value = 4297
$ws.Range("I126").Value = 4168.5713
$ws.Range("J126").Value = 4596.6665
$ws.Range("K126").Value = 12505.7139
$ws.Range("L126").Value = 13789.9995
$ws.Range("M126").Value = -10035.7139
$ws.Range("N126").Value = -18729.9995

$ws.Range("H132").Value = 6105027
$ws.Range("I132").Value = 6764349
$ws.Range("K132").Value = 20293047
$ws.Range("M132").Value = -20290517

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1493.6
$ws.Range("I96").Value = 1456.3334
$ws.Range("J96").Value = 1549.5
$ws.Range("K96").Value = 1456.3334
$ws.Range("L96").Value = 1549.5
$ws.Range("M96").Value = -83.33339999999998
$ws.Range("N96").Value = -4295.5

$ws.Range("H113").Value = 682.9091
$ws.Range("I113").Value = 643.1070999999999
$ws.Range("K113").Value = 1929.3213
$ws.Range("M113").Value = 240.6787000000002

$ws.Range("H132").Value = 9618551
$ws.Range("I132").Value = 11906796
$ws.Range("K132").Value = 35720388
$ws.Range("M132").Value = -35720388
